$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Cells.Item(34, 1).Value = 111635413
$ws.Cells.Item(34, 2).Value = 89369
$ws.Cells.Item(34, 4).Value = 'LC'
$ws.Cells.Item(34, 5).Value = 5447
$ws.Cells.Item(34, 6).Value = 'Vedticka'
$ws.Cells.Item(34, 7).Value = 'Fuscoporia viticola'
$ws.Cells.Item(34, 8).Value = '(Schwein.) Murrill'
$ws.Cells.Item(34, 17).Value = 539850.8116781802
$ws.Cells.Item(34, 18).Value = 7198361.834730743
$ws.Cells.Item(34, 26).Value = '09:57'
$ws.Cells.Item(34, 28).Value = '09:57'
$ws.Cells.Item(34, 49).Value = 'Yasmine Kindlund'
$ws.Cells.Item(34, 50).Value = 'Yasmine Kindlund, Isak Vahlström'

# Row 36
$ws.Cells.Item(36, 1).Value = 111634865
$ws.Cells.Item(36, 2).Value = 90087
$ws.Cells.Item(36, 4).Value = 'LC'
$ws.Cells.Item(36, 5).Value = 3298
$ws.Cells.Item(36, 6).Value = 'Trådticka'
$ws.Cells.Item(36, 7).Value = 'Climacocystis borealis'
$ws.Cells.Item(36, 8).Value = '(Fr.) Kotl. & Pouzar'
$ws.Cells.Item(36, 17).Value = 539879.8909062841
$ws.Cells.Item(36, 18).Value = 7198349.058794393
$ws.Cells.Item(36, 26).Value = '09:47'
$ws.Cells.Item(36, 28).Value = '09:47'
$ws.Cells.Item(36, 49).Value = 'Isak Vahlström'
$ws.Cells.Item(36, 50).Value = 'Isak Vahlström, Yasmine Kindlund'

# Row 37
$ws.Cells.Item(37, 1).Value = 111635461
$ws.Cells.Item(37, 2).Value = 89590
$ws.Cells.Item(37, 4).Value = 'VU'
$ws.Cells.Item(37, 5).Value = 48
$ws.Cells.Item(37, 6).Value = 'Lappticka'
$ws.Cells.Item(37, 7).Value = 'Amylocystis lapponica'
$ws.Cells.Item(37, 8).Value = '(Romell) Singer'
$ws.Cells.Item(37, 17).Value = 539846.9353019162
$ws.Cells.Item(37, 18).Value = 7198365.604689348
$ws.Cells.Item(37, 26).Value = '09:56'
$ws.Cells.Item(37, 28).Value = '09:56'
$ws.Cells.Item(37, 49).Value = 'Yasmine Kindlund'
$ws.Cells.Item(37, 50).Value = 'Yasmine Kindlund, Isak Vahlström'

# Row 38
$ws.Cells.Item(38, 1).Value = 111635419
$ws.Cells.Item(38, 2).Value = 89405
$ws.Cells.Item(38, 4).Value = 'NT'
$ws.Cells.Item(38, 5).Value = 1202
$ws.Cells.Item(38, 6).Value = 'Ullticka'
$ws.Cells.Item(38, 7).Value = 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(38, 8).Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(38, 17).Value = 539844.8100177459
$ws.Cells.Item(38, 18).Value = 7198365.57640036
$ws.Cells.Item(38, 26).Value = '09:58'
$ws.Cells.Item(38, 28).Value = '09:58'
$ws.Cells.Item(38, 49).Value = 'Yasmine Kindlund'
$ws.Cells.Item(38, 50).Value = 'Yasmine Kindlund, Isak Vahlström'

# Row 39
$ws.Cells.Item(39, 1).Value = 111634867
$ws.Cells.Item(39, 2).Value = 73696
$ws.Cells.Item(39, 4).Value = 'NT'
$ws.Cells.Item(39, 5).Value = 6440
$ws.Cells.Item(39, 6).Value = 'Vitgrynig nållav'
$ws.Cells.Item(39, 7).Value = 'Chaenotheca subroscida'
$ws.Cells.Item(39, 8).Value = '(Eitner) Zahlbr.'
$ws.Cells.Item(39, 17).Value = 539871.8034722162
$ws.Cells.Item(39, 18).Value = 7198349.800304586
$ws.Cells.Item(39, 26).Value = '09:46'
$ws.Cells.Item(39, 28).Value = '09:46'
$ws.Cells.Item(39, 49).Value = 'Isak Vahlström'
$ws.Cells.Item(39, 50).Value = 'Isak Vahlström, Yasmine Kindlund'

# Row 40
$ws.Cells.Item(40, 1).Value = 111634869
$ws.Cells.Item(40, 2).Value = 78578
$ws.Cells.Item(40, 4).Value = 'NT'
$ws.Cells.Item(40, 5).Value = 6458
$ws.Cells.Item(40, 6).Value = 'Lunglav'
$ws.Cells.Item(40, 7).Value = 'Lobaria pulmonaria'
$ws.Cells.Item(40, 8).Value = '(L.) Hoffm.'
$ws.Cells.Item(40, 17).Value = 539972.5933666634
$ws.Cells.Item(40, 18).Value = 7198379.169240371
$ws.Cells.Item(40, 26).Value = '09:27'
$ws.Cells.Item(40, 28).Value = '09:27'
$ws.Cells.Item(40, 49).Value = 'Isak Vahlström'
$ws.Cells.Item(40, 50).Value = 'Isak Vahlström, Yasmine Kindlund'

# Row 41
$ws.Cells.Item(41, 1).Value = 111635462
$ws.Cells.Item(41, 2).Value = 89590
$ws.Cells.Item(41, 4).Value = 'VU'
$ws.Cells.Item(41, 5).Value = 48
$ws.Cells.Item(41, 6).Value = 'Lappticka'
$ws.Cells.Item(41, 7).Value = 'Amylocystis lapponica'
$ws.Cells.Item(41, 8).Value = '(Romell) Singer'
$ws.Cells.Item(41, 17).Value = 539961.7289606878
$ws.Cells.Item(41, 18).Value = 7198365.011824355
$ws.Cells.Item(41, 26).Value = '09:30'
$ws.Cells.Item(41, 28).Value = '09:30'
$ws.Cells.Item(41, 49).Value = 'Yasmine Kindlund'
$ws.Cells.Item(41, 50).Value = 'Yasmine Kindlund, Isak Vahlström'

# Row 42
$ws.Cells.Item(42, 1).Value = 111634866
$ws.Cells.Item(42, 2).Value = 77515
$ws.Cells.Item(42, 4).Value = 'NT'
$ws.Cells.Item(42, 5).Value = 6425
$ws.Cells.Item(42, 6).Value = 'Garnlav'
$ws.Cells.Item(42, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(42, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(42, 17).Value = 539873.9909718054
$ws.Cells.Item(42, 18).Value = 7198345.158619706
$ws.Cells.Item(42, 26).Value = '09:46'
$ws.Cells.Item(42, 28).Value = '09:46'
$ws.Cells.Item(42, 49).Value = 'Isak Vahlström'
$ws.Cells.Item(42, 50).Value = 'Isak Vahlström, Yasmine Kindlund'

# Row 43
$ws.Cells.Item(43, 1).Value = 111635437
$ws.Cells.Item(43, 2).Value = 89845
$ws.Cells.Item(43, 4).Value = 'VU'
$ws.Cells.Item(43, 5).Value = 1209
$ws.Cells.Item(43, 6).Value = 'Rynkskinn'
$ws.Cells.Item(43, 7).Value = 'Phlebia centrifuga'
$ws.Cells.Item(43, 8).Value = 'P.Karst.'
$ws.Cells.Item(43, 17).Value = 539973.573864806
$ws.Cells.Item(43, 18).Value = 7198369.416147546
$ws.Cells.Item(43, 26).Value = '09:32'
$ws.Cells.Item(43, 28).Value = '09:32'
$ws.Cells.Item(43, 49).Value = 'Yasmine Kindlund'
$ws.Cells.Item(43, 50).Value = 'Yasmine Kindlund, Isak Vahlström'

# Row 44
$ws.Cells.Item(44, 1).Value = 111635452
$ws.Cells.Item(44, 2).Value = 78578
$ws.Cells.Item(44, 4).Value = 'NT'
$ws.Cells.Item(44, 5).Value = 6458
$ws.Cells.Item(44, 6).Value = 'Lunglav'
$ws.Cells.Item(44, 7).Value = 'Lobaria pulmonaria'
$ws.Cells.Item(44, 8).Value = '(L.) Hoffm.'
$ws.Cells.Item(44, 17).Value = 539861.2921981018
$ws.Cells.Item(44, 18).Value = 7198404.860384831
$ws.Cells.Item(44, 26).Value = '09:50'
$ws.Cells.Item(44, 28).Value = '09:50'
$ws.Cells.Item(44, 49).Value = 'Yasmine Kindlund'
$ws.Cells.Item(44, 50).Value = 'Yasmine Kindlund, Isak Vahlström'

# Row 45
$ws.Cells.Item(45, 1).Value = 111635499
$ws.Cells.Item(45, 2).Value = 85715
$ws.Cells.Item(45, 4).Value = 'NT'
$ws.Cells.Item(45, 5).Value = 510
$ws.Cells.Item(45, 6).Value = 'Doftskinn'
$ws.Cells.Item(45, 7).Value = 'Cystostereum murrayi'
$ws.Cells.Item(45, 8).Value = '(Berk. & M.A. Curtis.) Pouzar'
$ws.Cells.Item(45, 17).Value = 540009.9192712342
$ws.Cells.Item(45, 18).Value = 7198353.766191677
$ws.Cells.Item(45, 26).Value = '09:32'
$ws.Cells.Item(45, 28).Value = '09:32'
$ws.Cells.Item(45, 49).Value = 'Yasmine Kindlund'
$ws.Cells.Item(45, 50).Value = 'Yasmine Kindlund, Isak Vahlström'

# Row 47
$ws.Cells.Item(47, 1).Value = 111634859
$ws.Cells.Item(47, 2).Value = 77515
$ws.Cells.Item(47, 4).Value = 'NT'
$ws.Cells.Item(47, 5).Value = 6425
$ws.Cells.Item(47, 6).Value = 'Garnlav'
$ws.Cells.Item(47, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(47, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(47, 17).Value = 539847.161346367
$ws.Cells.Item(47, 18).Value = 7198348.622951495
$ws.Cells.Item(47, 26).Value = '09:58'
$ws.Cells.Item(47, 28).Value = '09:58'
$ws.Cells.Item(47, 49).Value = 'Isak Vahlström'
$ws.Cells.Item(47, 50).Value = 'Isak Vahlström, Yasmine Kindlund'

# Row 48
$ws.Cells.Item(48, 1).Value = 111635489
$ws.Cells.Item(48, 2).Value = 77515
$ws.Cells.Item(48, 4).Value = 'NT'
$ws.Cells.Item(48, 5).Value = 6425
$ws.Cells.Item(48, 6).Value = 'Garnlav'
$ws.Cells.Item(48, 7).Value = 'Alectoria sarmentosa'
$ws.Cells.Item(48, 8).Value = '(Ach.) Ach.'
$ws.Cells.Item(48, 17).Value = 539945.9506927577
$ws.Cells.Item(48, 18).Value = 7198336.776317291
$ws.Cells.Item(48, 26).Value = '09:39'
$ws.Cells.Item(48, 28).Value = '09:39'
$ws.Cells.Item(48, 49).Value = 'Yasmine Kindlund'
$ws.Cells.Item(48, 50).Value = 'Yasmine Kindlund, Isak Vahlström'

# Row 49
$ws.Cells.Item(49, 1).Value = 111635445
$ws.Cells.Item(49, 2).Value = 89686
$ws.Cells.Item(49, 4).Value = 'NT'
$ws.Cells.Item(49, 5).Value = 658
$ws.Cells.Item(49, 6).Value = 'Rosenticka'
$ws.Cells.Item(49, 7).Value = 'Rhodofomes roseus'
$ws.Cells.Item(49, 8).Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Cells.Item(49, 17).Value = 539972.1173992374
$ws.Cells.Item(49, 18).Value = 7198351.138093079
$ws.Cells.Item(49, 26).Value = '09:30'
$ws.Cells.Item(49, 28).Value = '09:30'
$ws.Cells.Item(49, 49).Value = 'Yasmine Kindlund'
$ws.Cells.Item(49, 50).Value = 'Yasmine Kindlund, Isak Vahlström'

# Row 50
$ws.Cells.Item(50, 1).Value = 111634868
$ws.Cells.Item(50, 2).Value = 78612
$ws.Cells.Item(50, 4).Value = 'LC'
$ws.Cells.Item(50, 5).Value = 6464
$ws.Cells.Item(50, 6).Value = 'Luddlav'
$ws.Cells.Item(50, 7).Value = 'Nephroma resupinatum'
$ws.Cells.Item(50, 8).Value = '(L.) Ach.'
$ws.Cells.Item(50, 17).Value = 539976.4302002029
$ws.Cells.Item(50, 18).Value = 7198378.371244119
$ws.Cells.Item(50, 26).Value = '09:28'
$ws.Cells.Item(50, 28).Value = '09:28'
$ws.Cells.Item(50, 49).Value = 'Isak Vahlström'
$ws.Cells.Item(50, 50).Value = 'Isak Vahlström, Yasmine Kindlund'
